# Added data to sample csv test:
# fill in the remaining header cells for the Weather/Play dataset
# (Outlook, Temperature already existed in A1:B1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Humidity"
$ws.Range("D1").Value = "Windy"
$ws.Range("E1").Value = "Play"

# Best-fit column B (holds "Temperature") like Excel does automatically
# when a header is the widest entry in its column.
$ws.Columns.Item(2).AutoFit()

# Move the selection off the header row, onto the first data row.
$ws.Range("A2").Select()
